# "Added Games Con info"
#
# Adds a new bulleted/numbered list item ("New article posted on surviving
# games con in Cologne, Germany") right after the existing "...game, Fe"
# list item, reusing that paragraph's list formatting (pStyle
# "ListParagraph", numPr ilvl/numId, Times New Roman 12pt run formatting).
# The document's hidden "_GoBack" bookmark - which sits collapsed right
# after the last character of the "...game, Fe" paragraph - needs to end
# up collapsed right after the last character of the newly added paragraph.

$d = $word.ActiveDocument

# The "_GoBack" bookmark marks the exact spot (collapsed range) where the
# new content belongs: right after "...game, Fe" and before its paragraph
# mark. Insert the new paragraph's text there, preceded by a unique marker
# that we then turn into a paragraph break via Find/Replace. Typing plain
# text directly at a collapsed bookmark's position naturally carries the
# bookmark forward so it stays attached to the end of the inserted text,
# which is exactly the placement the target document needs; turning the
# marker into "^p" afterwards splits that text into its own paragraph
# (inheriting the source paragraph's list formatting) without disturbing
# that bookmark tracking.
$bookmark = $d.Bookmarks.Item("_GoBack")
$insertionPoint = $d.Range($bookmark.Start, $bookmark.Start)
$insertionPoint.InsertAfter("~~SPLIT~~New article posted on surviving games con in Cologne, Germany")

$d.Content.Find.Execute("~~SPLIT~~", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "^p", 2)
